$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C for "Квартира" - shifts existing C..H to D..I
$ws.Range("C1:C3").EntireColumn.Insert()

# Insert three new columns after the (now shifted) "Показание" column (G)
# to hold tariff 2, 3, 4 values - shifts existing G..I (Житель, Источник) to J..K
$ws.Range("H1:J3").EntireColumn.Insert()

# Fill header row (row 1)
$ws.Range("C1").Value = "Квартира"
$ws.Range("G1").Value = "Показание по тарифу №1"
$ws.Range("H1").Value = "Показание по тарифу №2"
$ws.Range("I1").Value = "Показание по тарифу №3"
$ws.Range("J1").Value = "Показание по тарифу №4"

# Fill data row 2
$ws.Range("C2").Value = "{d.meter[i].unitName}"
$ws.Range("H2").Value = "{d.meter[i].value2}"
$ws.Range("I2").Value = "{d.meter[i].value3}"
$ws.Range("J2").Value = "{d.meter[i].value4}"

# Fill data row 3
$ws.Range("C3").Value = "{d.meter[i + 1].unitName}"
$ws.Range("H3").Value = "{d.meter[i + 1].value2}"
$ws.Range("I3").Value = "{d.meter[i + 1].value3}"
$ws.Range("J3").Value = "{d.meter[i + 1].value4}"
